$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 106-112 are being re-ordered / updated: three previously "pending"
# rows (F column = "a" / "r", no H value) were processed (F -> 1, H filled
# in), the now fully-processed "Bar Collar" row moved up, and the last
# pending row became "Leg Extension Machine" (still F = "r", no H value).

# Row 106: Ankle Straps (was row 109, now processed)
$ws.Range("A106").Value = 7
$ws.Range("B106").Value = "Ankle Straps"
$ws.Range("D106").Value = 5
$ws.Range("F106").Value = 1
$ws.Range("H106").Value = 117

# Row 107: Bar Collar (was row 106)
$ws.Range("A107").Value = 14
$ws.Range("B107").Value = "Bar Collar"
$ws.Range("D107").Value = 11
$ws.Range("F107").Value = 1
$ws.Range("H107").Value = 127

# Row 108: Duck Walk Handle (was row 107)
$ws.Range("A108").Value = 39
$ws.Range("B108").Value = "Duck Walk Handle"
$ws.Range("D108").Value = 35
$ws.Range("F108").Value = 1
$ws.Range("H108").Value = 50

# Row 109: Neck machine (was row 110, now processed)
$ws.Range("A109").Value = 67
$ws.Range("B109").Value = "Neck machine"
$ws.Range("D109").Value = 62
$ws.Range("F109").Value = 1
$ws.Range("H109").Value = 117

# Row 110: Resistance Band (was row 108)
$ws.Range("A110").Value = 80
$ws.Range("B110").Value = "Resistance Band"
$ws.Range("D110").Value = 75
$ws.Range("F110").Value = 1
$ws.Range("H110").Value = 134

# Row 111: Twist Machine (was row 112, now processed)
$ws.Range("A111").Value = 109
$ws.Range("B111").Value = "Twist Machine"
$ws.Range("D111").Value = 104
$ws.Range("F111").Value = 1
$ws.Range("H111").Value = 81

# Row 112: Leg Extension Machine (was row 111, still pending/"r")
$ws.Range("A112").Value = 61
$ws.Range("B112").Value = "Leg Extension Machine"
$ws.Range("D112").Value = 57
$ws.Range("F112").Value = "r"
$ws.Range("H112").ClearContents()
